# Regenerate merged AHB files
#
# 1. Rename the header row's shared strings:
#      *_old -> *_FV2304
#      *_new -> *_FV2310   (note: "Bedinung_new" -> "Bedinung_FV2310", keeping the
#                            original "Bedinung" spelling, not "Bedingung")
# 2. Turn the used range A1:T72 into a native Excel Table ("Table1").
# 3. Freeze the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- 1. Rename header cells (A1:T1) ---------------------------------------
$headers = @(
  "Segmentname_FV2304",
  "Segmentgruppe_FV2304",
  "Segment_FV2304",
  "Datenelement_FV2304",
  "Segment ID_FV2304",
  "Code_FV2304",
  "Qualifier_FV2304",
  "Beschreibung_FV2304",
  "Bedingungsausdruck_FV2304",
  "Bedingung_FV2304",
  "diff",
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedinung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Convert the used range into a native table -------------------------
$tableRange = $ws.Range("A1:T72")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Headers renamed, Table1 added, header row frozen."
